# CI Build: Synchronisation of the project's models with the google sheet metadata version 9841
# Improve unit cost handling
#
# For every sheet that has a "unitCost" column followed by "qualityGrading",
# insert two new columns right after "unitCost":
#   - unitCostCurrency
#   - unitCostNote

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Service",
    "Product",
    "Antibody",
    "Hybridoma",
    "Protein",
    "NucleicAcid",
    "DetectionKit",
    "Bundle",
    "Virus",
    "Bacterium",
    "Fungus",
    "Protozoan",
    "Viroid",
    "Prion"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $found = $ws.Rows.Item(1).Find("unitCost")
    $unitCostCol = $found.Column
    $insertCol = $unitCostCol + 1

    # Insert two blank columns right after the "unitCost" column
    $ws.Cells.Item(1, $insertCol).EntireColumn.Insert()
    $ws.Cells.Item(1, $insertCol).EntireColumn.Insert()

    $ws.Cells.Item(1, $insertCol).Value2 = "unitCostCurrency"
    $ws.Cells.Item(1, $insertCol + 1).Value2 = "unitCostNote"
}
